$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.23349404335022
$ws.Range("B1").Value = 2.294487476348877
$ws.Range("C1").Value = 3.463886260986328
$ws.Range("D1").Value = 3.883548498153687
$ws.Range("E1").Value = 1.04381000995636
